$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells involved in this update to remain plain text,
# so numeric-looking strings (e.g. "1.00", "305.80") are not coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.667.19"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.279.64"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.80"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.44"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.506"
$ws.Range("E7").Value = "  -2.76%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").Value = "  -3.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.47"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.22"
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.631.46"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.278.44"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.778"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.581.85"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.91"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.07"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.13"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.04"
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.15"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.61"
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0690"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.999.33"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0279"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.22"
$ws.Range("E45").Value = "  +4.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.96"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("E47").Value = "  -7.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.77"
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.87"
$ws.Range("E49").Value = "  +5.07%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.62"
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.500.84"
$ws.Range("E51").Value = "  -1.20%  "
